$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values for rows 2-4
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Update column B value for row 4
$ws.Range("B4").Value = 2

# Delete row 5 entirely (shifts cells up), removing A5/B5 values
$ws.Rows.Item(5).Delete()
